$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 1866.6666
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 1866.6666
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H137").Value = 5750
$ws.Range("I137").Value = 10000
$ws.Range("J137").Value = 1500
$ws.Range("K137").Value = 30000
$ws.Range("L137").Value = 4500
$ws.Range("M137").Value = -27450
$ws.Range("N137").Value = -9600
$ws.Range("H138").Value = 2972.52
$ws.Range("I138").Value = 1361.7142
$ws.Range("J138").Value = 3400.7087
$ws.Range("K138").Value = 4085.1426
$ws.Range("L138").Value = 10202.1261
$ws.Range("M138").Value = 1054.8574
$ws.Range("N138").Value = -20482.1261
$ws.Range("H141").Value = 4877.778
$ws.Range("I141").Value = 4816.6665
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 14449.9995
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -9269.999500000002
$ws.Range("N141").Value = -25360

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8054.391
$ws.Range("I32").Value = 5979.3037
$ws.Range("J32").Value = 16993.23
$ws.Range("K32").Value = 5979.3037
$ws.Range("L32").Value = 16993.23
$ws.Range("M32").Value = -5692.3037
$ws.Range("N32").Value = -17567.23
$ws.Range("H74").Value = 925.31744
$ws.Range("I74").Value = 887.0179000000001
$ws.Range("J74").Value = 1231.7142
$ws.Range("K74").Value = 887.0179000000001
$ws.Range("L74").Value = 1231.7142
$ws.Range("M74").Value = -13.01790000000005
$ws.Range("N74").Value = -2979.7142
$ws.Range("H77").Value = 925.31744
$ws.Range("I77").Value = 887.0179000000001
$ws.Range("J77").Value = 1231.7142
$ws.Range("K77").Value = 4435.0895
$ws.Range("L77").Value = 6158.571
$ws.Range("M77").Value = -67.08950000000004
$ws.Range("N77").Value = -14894.571
$ws.Range("H97").Value = 984.2
$ws.Range("I97").Value = 702.5
$ws.Range("J97").Value = 2111
$ws.Range("K97").Value = 702.5
$ws.Range("L97").Value = 2111
$ws.Range("M97").Value = -206.5
$ws.Range("N97").Value = -3103
$ws.Range("H132").Value = 2136.6
$ws.Range("I132").Value = 1518.1666
$ws.Range("J132").Value = 3726.8572
$ws.Range("K132").Value = 4554.4998
$ws.Range("L132").Value = 11180.5716
$ws.Range("M132").Value = -2024.4998
$ws.Range("N132").Value = -16240.5716

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 762.63635
$ws.Range("I94").Value = 843.7778
$ws.Range("K94").Value = 843.7778
$ws.Range("M94").Value = -392.7778

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H31").Value = 2579.3958
$ws.Range("I31").Value = 2596.6511
$ws.Range("J31").Value = 2431
$ws.Range("K31").Value = 2596.6511
$ws.Range("L31").Value = 2431
$ws.Range("M31").Value = -2301.6511
$ws.Range("N31").Value = -3021
$ws.Range("H34").Value = 2579.3958
$ws.Range("I34").Value = 2596.6511
$ws.Range("J34").Value = 2431
$ws.Range("K34").Value = 2596.6511
$ws.Range("L34").Value = 2431
$ws.Range("M34").Value = -2394.6511
$ws.Range("N34").Value = -2835
$ws.Range("H38").Value = 12000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 12000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 12000
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -12754
$ws.Range("H46").Value = 12000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 12000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 12000
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -12422
$ws.Range("H132").Value = 1542.6666
$ws.Range("I132").Value = 1160.3334
$ws.Range("J132").Value = 2307.3333
$ws.Range("K132").Value = 3481.0002
$ws.Range("L132").Value = 6921.999899999999
$ws.Range("M132").Value = -951.0001999999999
$ws.Range("N132").Value = -11981.9999
$ws.Range("H134").Value = 15152599
$ws.Range("I134").Value = 1040
$ws.Range("J134").Value = 62501220
$ws.Range("K134").Value = 3120
$ws.Range("L134").Value = 187503660
$ws.Range("M134").Value = -585
$ws.Range("N134").Value = -187508730

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 18556998
$ws.Range("J131").Value = 1421.0526
$ws.Range("L131").Value = 4263.1578
$ws.Range("N131").Value = -14343.1578

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3798.762
$ws.Range("I80").Value = 3702.3076
$ws.Range("J80").Value = 3955.5
$ws.Range("K80").Value = 3702.3076
$ws.Range("L80").Value = 3955.5
$ws.Range("M80").Value = -2704.3076
$ws.Range("N80").Value = -5951.5
$ws.Range("H83").Value = 3798.762
$ws.Range("I83").Value = 3702.3076
$ws.Range("J83").Value = 3955.5
$ws.Range("K83").Value = 18511.538
$ws.Range("L83").Value = 19777.5
$ws.Range("M83").Value = -13519.538
$ws.Range("N83").Value = -29761.5
$ws.Range("H132").Value = 5045.364
$ws.Range("I132").Value = 5500.1665
$ws.Range("J132").Value = 4499.6
$ws.Range("K132").Value = 16500.4995
$ws.Range("L132").Value = 13498.8
$ws.Range("M132").Value = -13970.4995
$ws.Range("N132").Value = -18558.8

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1511.8422
$ws.Range("J22").Value = 1573.6111
$ws.Range("L22").Value = 1573.6111
$ws.Range("N22").Value = -2163.6111
$ws.Range("H27").Value = 1511.8422
$ws.Range("J27").Value = 1573.6111
$ws.Range("L27").Value = 1573.6111
$ws.Range("N27").Value = -1787.6111
$ws.Range("H46").Value = 15386018
$ws.Range("I46").Value = 18183136
$ws.Range("J46").Value = 1875
$ws.Range("K46").Value = 18183136
$ws.Range("L46").Value = 1875
$ws.Range("M46").Value = -18182948
$ws.Range("N46").Value = -2251

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 44005.5
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 44005.5
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 44005.5
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -44485.5
$ws.Range("H81").Value = 1938
$ws.Range("I81").Value = 1672.5
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 3345
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -2284
$ws.Range("N81").Value = -8122
$ws.Range("H84").Value = 1938
$ws.Range("I84").Value = 1672.5
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 16725
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -11421
$ws.Range("N84").Value = -40608
